# Automatische test-sync: 2025-08-04 20:16:50
# Adds a new "Inkoop / Bestellingen" test-mail row to the Logs sheet and
# the corresponding rollup row to the Dashboard sheet, extends the
# conditional-formatting ranges to cover the new row, and widens the
# chart's category/value source ranges so the new row is plotted too.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 10 -------------------------------------------
$logs.Cells.Item(10, 1).Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Cells.Item(10, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(10, 3).Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$logs.Cells.Item(10, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item(10, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item(10, 6).Value = "2025-08-04 20:16:21"
$logs.Cells.Item(10, 7).Value = "Ja"
$logs.Cells.Item(10, 8).Value = "Ja"
$logs.Cells.Item(10, 9).Value = "Nee"
$logs.Cells.Item(10, 10).Value = "Nee"

# --- Dashboard sheet: append rollup row 5 --------------------------------
$dash.Cells.Item(5, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(5, 2).Value = 1

# --- Widen the conditional-formatting ranges to include the new row -----
# (ModifyAppliesToRange keeps the existing rules - priorities, dxfIds,
# formulas - intact and only changes the sqref, unlike Delete+Add.)
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $col + "2:" + $col + "9"
    $newRange = $col + "2:" + $col + "10"
    $rules = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Extend the chart's category/value series to include the new row ----
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"
